$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format before assignment, otherwise Excel auto-converts them to numeric
# values (losing trailing zeros / exact formatting). We then restore the
# cell style back to Normal so no stray style/number-format diff is left behind.
$numericCells = @(
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D17',
    'D19',
    'D20',
    'D21',
    'D23',
    'D24',
    'D25',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D38',
    'D40',
    'D42',
    'D43',
    'D44',
)
foreach ($cell in $numericCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply all the updated values
$ws.Range('D2').Value = '52.072.24'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '2.819.36'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '355.73'
$ws.Range('E5').Value = '  +2.83%  '
$ws.Range('D6').Value = '111.98'
$ws.Range('E6').Value = '  -4.13%  '
$ws.Range('D7').Value = '0.571'
$ws.Range('E7').Value = '  +3.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('D10').Value = '40.81'
$ws.Range('E10').Value = '  -5.90%  '
$ws.Range('D11').Value = '0.0861'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '19.93'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = '7.76'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '3.260.50'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '2.827.02'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '0.932'
$ws.Range('E17').Value = '  +4.44%  '
$ws.Range('D18').Value = '51.850.32'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = '7.52'
$ws.Range('E19').Value = '  +5.51%  '
$ws.Range('D20').Value = '3.20'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = '13.45'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = '0.0₃0996'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').Value = '70.83'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = '269.36'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.144'
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('B31').Value = 'VeChain'
$ws.Range('C31').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D31').Value = '0.0488'
$ws.Range('E31').Value = '  +17.71%  '
$ws.Range('D32').Value = '52.42'
$ws.Range('E32').Value = '  +4.23%  '
$ws.Range('D33').Value = '34.81'
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('D34').Value = '5.94'
$ws.Range('E34').Value = '  +3.78%  '
$ws.Range('D35').Value = '5.61'
$ws.Range('E35').Value = '  +12.48%  '
$ws.Range('D36').Value = '0.0854'
$ws.Range('E36').Value = '  +3.65%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').Value = '3.28'
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  -4.24%  '
$ws.Range('D40').Value = '18.36'
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('D42').Value = '127.29'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').Value = '23.25'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '2.49'
$ws.Range('E44').Value = '  -7.49%  '
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = '2.080.06'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  -5.83%  '
$ws.Range('E49').Value = '  +5.92%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('E51').Value = '  +2.02%  '

# Restore default styling on the cells we temporarily reformatted as Text
foreach ($cell in $numericCells) {
    $ws.Range($cell).Style = "Normal"
}
